$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Content fix: rename the "CDM" construction-date bucket to "CDL"
# across the three long mapping-scheme text cells (Offices/Trade/Hotels).
foreach ($addr in @("B2", "C2", "D2")) {
    $cell = $ws.Range($addr)
    $text = $cell.Value2
    $cell.Value2 = $text.Replace("CDM", "CDL")
}

# --- Layout: widen the three data columns so the long strings are readable
$ws.Columns.Item(2).ColumnWidth = 40.666666666666664   # -> width 41.5
$ws.Columns.Item(3).ColumnWidth = 35.5                 # -> width ~36.332
$ws.Columns.Item(4).ColumnWidth = 33                   # -> width ~33.832

# --- Layout: grow row 2 to the maximum height and wrap the long text
$ws.Rows.Item(2).RowHeight = 409.6
$ws.Range("B2:D2").WrapText = $true

# --- Selection moves to B5 as left by the editor
$ws.Range("B5").Select() | Out-Null
